$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Row 2: two more schedule dates (J2 = 17/06/2020, K2 = 18/06/2020) ---
# Values first, then copy the date-column format (from I2) so the same
# existing style ("s=4") is reused instead of minting a new one.
$ws.Range("J2").Value = 43999
$ws.Range("K2").Value = 44000
$ws.Range("I2").Copy()
$ws.Range("J2:K2").PasteSpecial(-4122) # xlPasteFormats

# --- Row 11: mark the two new days as "in progress" with the yellow fill ---
# used elsewhere for task markers, but without the extra alignment override
# (this mints the one new cellXf the workbook needs).
$ws.Range("J11").Interior.ColorIndex = 6
$ws.Range("K11").Interior.ColorIndex = 6

# --- Row 12: hours logged for the two new days ---
$ws.Range("I12").Copy()
$ws.Range("J12").PasteSpecial(-4122)
$ws.Range("J12").Value = "4.5 h."

$ws.Range("I12").Copy()
$ws.Range("K12").PasteSpecial(-4122)
$ws.Range("K12").Value = "2.5 h."

$excel.CutCopyMode = 0

# --- Restore the recorded selection/view state ---
$ws.Range("I20").Select()

$wb.Save()
